$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2381042689085007
$ws.Range("B1").Value = 0.2103091180324554
$ws.Range("C1").Value = 0.200173631310463
$ws.Range("D1").Value = 0.2375946491956711
$ws.Range("E1").Value = 0.3300251662731171
